$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "70.351.30"
$ws.Range("E2").Value2 = "  +0.84%  "
$ws.Range("D3").Value2 = "3.621.35"
$ws.Range("E3").Value2 = "  +2.76%  "
$ws.Range("E4").Value2 = "  +0.13%  "
$ws.Range("D5").Value2 = "603.09"
$ws.Range("E5").Value2 = "  -0.72%  "
$ws.Range("D6").Value2 = "196.02"
$ws.Range("E6").Value2 = "  -0.29%  "
$ws.Range("E7").Value2 = "  -0.91%  "
$ws.Range("E8").Value2 = "  +0.09%  "
$ws.Range("D9").Value2 = "0.215"
$ws.Range("E9").Value2 = "  +6.88%  "
$ws.Range("D10").Value2 = "0.643"
$ws.Range("E10").Value2 = "  -1.07%  "
$ws.Range("D11").Value2 = "53.23"
$ws.Range("E11").Value2 = "  -1.15%  "
$ws.Range("D12").Value2 = "0.0000305"
$ws.Range("E12").Value2 = "  +0.66%  "
$ws.Range("D13").Value2 = "9.52"
$ws.Range("E13").Value2 = "  -0.10%  "
$ws.Range("D14").Value2 = "4.197.28"
$ws.Range("E14").Value2 = "  +2.88%  "
$ws.Range("D15").Value2 = "606.64"
$ws.Range("E15").Value2 = "  +1.61%  "
$ws.Range("E16").Value2 = "  +0.97%  "
$ws.Range("D17").Value2 = "70.475.62"
$ws.Range("E17").Value2 = "  +0.81%  "
$ws.Range("B18").Value2 = "Chainlink"
$ws.Range("C18").Value2 = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value2 = "18.98"
$ws.Range("E18").Value2 = "  -0.81%  "
$ws.Range("B19").Value2 = "WrappedEther"
$ws.Range("C19").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value2 = "3.590.44"
$ws.Range("E19").Value2 = "  +1.99%  "
$ws.Range("E20").Value2 = "  +1.05%  "
$ws.Range("D21").Value2 = "0.996"
$ws.Range("E21").Value2 = "  +0.21%  "
$ws.Range("D22").Value2 = "18.12"
$ws.Range("E22").Value2 = "  -0.99%  "
$ws.Range("D23").Value2 = "5.28"
$ws.Range("E23").Value2 = "  -0.35%  "
$ws.Range("D24").Value2 = "103.58"
$ws.Range("E24").Value2 = "  +1.10%  "
$ws.Range("E25").Value2 = "  -1.52%  "
$ws.Range("E26").Value2 = "  -6.38%  "
$ws.Range("D27").Value2 = "10.59"
$ws.Range("E27").Value2 = "  -2.53%  "
$ws.Range("D28").Value2 = "9.70"
$ws.Range("E28").Value2 = "  +0.77%  "
$ws.Range("D29").Value2 = "33.75"
$ws.Range("E29").Value2 = "  +0.72%  "
$ws.Range("D30").Value2 = "4.68"
$ws.Range("E30").Value2 = "  +8.65%  "
$ws.Range("D31").Value2 = "7.18"
$ws.Range("E31").Value2 = "  +1.41%  "
$ws.Range("E32").Value2 = "  -1.67%  "
$ws.Range("E33").Value2 = "  -0.02%  "
$ws.Range("D34").Value2 = "63.27"
$ws.Range("E34").Value2 = "  +0.24%  "
$ws.Range("D35").Value2 = "0.0₃0887"
$ws.Range("E35").Value2 = "  +3.33%  "
$ws.Range("D36").Value2 = "3.959.79"
$ws.Range("E36").Value2 = "  +5.92%  "
$ws.Range("E37").Value2 = "  +0.21%  "
$ws.Range("D38").Value2 = "3.06"
$ws.Range("E38").Value2 = "  -0.49%  "
$ws.Range("D39").Value2 = "516.08"
$ws.Range("E39").Value2 = "  +5.84%  "
$ws.Range("D40").Value2 = "0.388"
$ws.Range("E40").Value2 = "  -1.29%  "
$ws.Range("D41").Value2 = "36.56"
$ws.Range("E41").Value2 = "  -0.20%  "
$ws.Range("E42").Value2 = "  -2.83%  "
$ws.Range("D43").Value2 = "0.136"
$ws.Range("E43").Value2 = "  +2.01%  "
$ws.Range("E44").Value2 = "  +1.25%  "
$ws.Range("D45").Value2 = "3.47"
$ws.Range("E45").Value2 = "  +5.59%  "
$ws.Range("D46").Value2 = "2.89"
$ws.Range("E46").Value2 = "  +2.57%  "
$ws.Range("E47").Value2 = "  +0.00%  "
$ws.Range("D48").Value2 = "8.55"
$ws.Range("E48").Value2 = "  +0.35%  "
$ws.Range("D50").Value2 = "0.000249"
$ws.Range("E50").Value2 = "  +0.87%  "
$ws.Range("E51").Value2 = "  +0.22%  "
